$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

# The "1a567e4d..." file's Status ("Ready for handoff") becomes "Handback transform failed"
# on every sheet that shows that status (Overview + the per-locale detail sheets all shared
# the same underlying string).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# zh-cn sheet: update status and record the error detail for the handback mismatch
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("K3").Value = "Handback file name: 2jv1zzty.sop is different with handoff file name: 1a567e4d-6b75-4fd1-be20-eb7581ce8c2e.de1afceb0086ef6f20ae439214eb25bfd3e68bfb.zh-cn."

# de-de sheet: update status and record the error detail for the handback mismatch
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("K3").Value = "Handback file name: 2jv1zzty.sop is different with handoff file name: 1a567e4d-6b75-4fd1-be20-eb7581ce8c2e.de1afceb0086ef6f20ae439214eb25bfd3e68bfb.de-de."
